$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04763786555579896
$ws.Range("C2").Value = 0.3127903958511391
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("G2").Value = 9.178186040332873

$ws.Range("B3").Value = 0.6753301551942219
$ws.Range("C3").Value = 114.8270160096505
$ws.Range("D3").Value = 26.21740644021617
$ws.Range("E3").Value = 645.3272768299601
$ws.Range("G3").Value = 787.0470294350209
